$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'isophonics_291'
$ws.Range("B2").Value = 'schubert-winterreise_114'
$ws.Range("C2").Value = 0.2392241379310345
$ws.Range("D2").Value = '[[''D'', ''G'', ''D'']]'
$ws.Range("E2").Value = '[[''D:maj/F#'', ''G:maj'', ''D:maj'']]'
$ws.Range("F2").Value = '[(5.20815, 11.03585)]'
$ws.Range("G2").Value = '[(57.48, 64.58)]'
$ws.Range("H2").Value = 'spotify:track:06ypiqmILMdVeaiErMFA91'
$ws.Range("I2").Value = ''

# Row 3
$ws.Range("A3").Value = 'schubert-winterreise_195'
$ws.Range("B3").Value = 'schubert-winterreise_186'
$ws.Range("C3").Value = 0.1714285714285714
$ws.Range("D3").Value = '[[''F:min/C'', ''C'', ''F:min/C'']]'
$ws.Range("E3").Value = '[[''F:min'', ''C:maj'', ''F:min'']]'
$ws.Range("F3").Value = '[(44.16, 48.16)]'
$ws.Range("G3").Value = '[(0.24, 5.08)]'
$ws.Range("H3").Value = ''
$ws.Range("I3").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

# Row 4
$ws.Range("A4").Value = 'isophonics_192'
$ws.Range("B4").Value = 'isophonics_81'
$ws.Range("C4").Value = 0.1644736842105263
$ws.Range("D4").Value = '[[''F'', ''G'', ''C'']]'
$ws.Range("E4").Value = '[[''E'', ''F#'', ''B'']]'
$ws.Range("F4").Value = '[(75.355124, 80.742154)]'
$ws.Range("G4").Value = '[(26.48873, 32.328548)]'
$ws.Range("H4").Value = ''
$ws.Range("I4").Value = ''

# Row 5
$ws.Range("A5").Value = 'isophonics_213'
$ws.Range("B5").Value = 'jaah_25'
$ws.Range("C5").Value = 0.1422413793103448
$ws.Range("D5").Value = '[[''E'', ''E:7'', ''A'']]'
$ws.Range("E5").Value = '[[''Bb'', ''Bb:7'', ''Eb'']]'
$ws.Range("F5").Value = '[(16.192174, 27.488682)]'
$ws.Range("G5").Value = '[(35.04, 40.71)]'
$ws.Range("H5").Value = ''
$ws.Range("I5").Value = ''

# Row 6
$ws.Range("A6").Value = 'schubert-winterreise_194'
$ws.Range("B6").Value = 'schubert-winterreise_93'
$ws.Range("C6").Value = 0.1712473572938689
$ws.Range("D6").Value = '[[''D:hdim7/F'', ''G:(3,5,b7,b9)'', ''C:min'']]'
$ws.Range("E6").Value = '[[''C:hdim7/D#'', ''F:(3,5,b7,b9)'', ''A#:min'']]'
$ws.Range("F6").Value = '[(21.76, 27.24)]'
$ws.Range("G6").Value = '[(60.42, 64.92)]'
$ws.Range("H6").Value = ''
$ws.Range("I6").Value = 'spotify:track:2qCvEz2hEb92VFATqVvrht'

# Row 7
$ws.Range("A7").Value = 'jaah_29'
$ws.Range("B7").Value = 'isophonics_135'
$ws.Range("C7").Value = 0.1366459627329192
$ws.Range("D7").Value = '[[''Ab'', ''Ab'', ''Db'']]'
$ws.Range("E7").Value = '[[''E'', ''E'', ''A'']]'
$ws.Range("F7").Value = '[(116.69, 119.66)]'
$ws.Range("G7").Value = '[(0.866546, 11.872804)]'
$ws.Range("H7").Value = ''
$ws.Range("I7").Value = ''

# Row 8
$ws.Range("A8").Value = 'isophonics_123'
$ws.Range("B8").Value = 'isophonics_194'
$ws.Range("C8").Value = 0.1230195712954334
$ws.Range("D8").Value = '[[''A'', ''E'', ''A'', ''E'']]'
$ws.Range("E8").Value = '[[''G'', ''D'', ''G'', ''D'']]'
$ws.Range("F8").Value = '[(9.156734, 13.057687)]'
$ws.Range("G8").Value = '[(142.379117, 149.240614)]'
$ws.Range("H8").Value = ''
$ws.Range("I8").Value = ''

# Row 9
$ws.Range("A9").Value = 'schubert-winterreise_28'
$ws.Range("B9").Value = 'schubert-winterreise_200'
$ws.Range("C9").Value = 0.5397727272727273
$ws.Range("D9").Value = '[[''A:maj/E'', ''E:7'', ''A:maj'', ''E:7'', ''A:maj'']]'
$ws.Range("E9").Value = '[[''E:maj'', ''B:7'', ''E:maj'', ''B:7'', ''E:maj'']]'
$ws.Range("F9").Value = '[(15.3, 19.72)]'
$ws.Range("G9").Value = '[(24.28, 44.36)]'
$ws.Range("H9").Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'
$ws.Range("I9").Value = ''

# Row 10
$ws.Range("A10").Value = 'isophonics_52'
$ws.Range("B10").Value = 'isophonics_87'
$ws.Range("C10").Value = 0.06166943521594684
$ws.Range("D10").Value = '[[''G'', ''G'', ''C'']]'
$ws.Range("E10").Value = '[[''A'', ''A'', ''D'']]'
$ws.Range("F10").Value = '[(0.497838, 5.89263)]'
$ws.Range("G10").Value = '[(59.508331, 64.685555)]'
$ws.Range("H10").Value = ''
$ws.Range("I10").Value = ''

# Row 11
$ws.Range("A11").Value = 'isophonics_277'
$ws.Range("B11").Value = 'isophonics_156'
$ws.Range("C11").Value = 0.1843137254901961
$ws.Range("D11").Value = '[[''B'', ''E'', ''E'', ''B'']]'
$ws.Range("E11").Value = '[[''E'', ''A'', ''A'', ''E'']]'
$ws.Range("F11").Value = '[(14.402743, 19.871043)]'
$ws.Range("G11").Value = '[(1.274255, 7.505864)]'
$ws.Range("H11").Value = 'spotify:track:2RnPATK99oGOZygnD2GTO6'
$ws.Range("I11").Value = ''

# Row 12
$ws.Range("A12").Value = 'isophonics_193'
$ws.Range("B12").Value = 'isophonics_275'
$ws.Range("C12").Value = 0.1792207792207792
$ws.Range("D12").Value = '[[''Ab'', ''Eb'', ''Bb/3''], [''Bb'', ''Eb'', ''Ab'']]'
$ws.Range("E12").Value = '[[''C'', ''G'', ''D''], [''D'', ''G'', ''C'']]'
$ws.Range("F12").Value = '[(25.942, 29.102), (10.286, 15.006)]'
$ws.Range("G12").Value = '[(15.298925, 24.401147), (28.800769, 40.713165)]'
$ws.Range("H12").Value = ''
$ws.Range("I12").Value = ''

# Row 13
$ws.Range("A13").Value = 'jaah_0'
$ws.Range("B13").Value = 'isophonics_56'
$ws.Range("C13").Value = 0.08233638282899367
$ws.Range("D13").Value = '[[''Eb:7'', ''Ab'', ''Ab''], [''Ab:maj6'', ''Bb:7'', ''Eb:7''], [''Bb:7'', ''Eb:7'', ''Ab'']]'
$ws.Range("E13").Value = '[[''G:7'', ''C'', ''C/7''], [''A:min7'', ''D:7'', ''G:7''], [''D:7'', ''G:7'', ''C'']]'
$ws.Range("F13").Value = '[(64.32, 68.27), (52.5, 56.45), (23.91, 25.87)]'
$ws.Range("G13").Value = '[(12.318208, 18.355396), (115.728276, 124.366099), (118.584331, 127.280204)]'
$ws.Range("H13").Value = ''
$ws.Range("I13").Value = ''

# Row 14
$ws.Range("A14").Value = 'schubert-winterreise_151'
$ws.Range("B14").Value = 'schubert-winterreise_74'
$ws.Range("C14").Value = 0.5333333333333333
$ws.Range("D14").Value = '[[''C:maj/G'', ''F:maj'', ''C:maj/G'', ''F:maj/A'']]'
$ws.Range("E14").Value = '[[''F:maj'', ''A#:maj'', ''F:maj'', ''A#:maj'']]'
$ws.Range("F14").Value = '[(117.54, 122.7)]'
$ws.Range("G14").Value = '[(129.38, 134.74)]'
$ws.Range("H14").Value = ''
$ws.Range("I14").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

# Row 15
$ws.Range("A15").Value = 'schubert-winterreise_11'
$ws.Range("B15").Value = 'schubert-winterreise_54'
$ws.Range("C15").Value = 0.2701149425287356
$ws.Range("D15").Value = '[[''C:min'', ''G:maj'', ''C:min'', ''G:maj'', ''C:min'', ''G:maj'']]'
$ws.Range("E15").Value = '[[''G:min'', ''D:maj/G'', ''G:min'', ''D:maj/G'', ''G:min'', ''D:maj/G'']]'
$ws.Range("F15").Value = '[(7.96, 21.56)]'
$ws.Range("G15").Value = '[(25.64, 55.14)]'
$ws.Range("H15").Value = ''
$ws.Range("I15").Value = ''

# Row 16
$ws.Range("A16").Value = 'isophonics_5'
$ws.Range("B16").Value = 'isophonics_112'
$ws.Range("C16").Value = 0.1875
$ws.Range("D16").Value = '[[''E'', ''A'', ''E'', ''A'', ''E'']]'
$ws.Range("E16").Value = '[[''A'', ''D/5'', ''A'', ''D/5'', ''A'']]'
$ws.Range("F16").Value = '[(9.162102, 17.463236)]'
$ws.Range("G16").Value = '[(0.421247, 7.703786)]'
$ws.Range("H16").Value = ''
$ws.Range("I16").Value = ''

# Row 17
$ws.Range("A17").Value = 'jaah_59'
$ws.Range("B17").Value = 'jaah_51'
$ws.Range("C17").Value = 0.02135157545605307
$ws.Range("D17").Value = '[[''C:7'', ''C:min7'', ''F:7'']]'
$ws.Range("E17").Value = '[[''D:7'', ''D:min7'', ''G:7'']]'
$ws.Range("F17").Value = '[(19.63, 22.36)]'
$ws.Range("G17").Value = '[(26.28, 30.35)]'
$ws.Range("H17").Value = ''
$ws.Range("I17").Value = ''
